$wb = $excel.ActiveWorkbook

# --- Sheet "СИП": mark the "Кол-во фидеров" header as required, like the
#     other mandatory columns (append a bold red "*" suffix) ---
$wsSip = $wb.Worksheets.Item("СИП")
[void]$wsSip.Activate()

$cell = $wsSip.Range("C1")
$cell.Value = "Кол-во фидеров*"
$star = $cell.Characters(15, 1)
$star.Font.Bold = $true
$star.Font.Size = 12
$star.Font.Color = 255

# Move the active selection on the "СИП" sheet
$null = $wsSip.Range("E4").Select()

# --- Sheet "Супервайзеры": move the active selection ---
$wsSuperv = $wb.Worksheets.Item("Супервайзеры")
[void]$wsSuperv.Activate()
$null = $wsSuperv.Range("E18").Select()

# Restore "СИП" as the active sheet/tab
[void]$wsSip.Activate()
